# Test Mail.xlsx update:
#  - Strip the bold/bordered/centred header style from row 1 (A1:D1) so the
#    header cells fall back to the default (unstyled) cell format.
#  - Refresh the two data rows with new test values (new sender addresses,
#    new "sent" timestamps/dates) produced by a later automation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header): remove the bold font + thin border + centered alignment
# that used to be applied to A1:D1 - cells go back to plain/default format.
$ws.Range("A1:D1").ClearFormats()

# Row 2: rajan@finlytyx.com -> adil@finlytyx.com, new send date/time.
$ws.Range("A2").Value = "adil@finlytyx.com"
$ws.Range("B2").Value = 46056
$ws.Range("B2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C2").Value = "Sent at 2026-02-03 17:09:30"

# Row 3: gokuldas@finlytyx.com -> mashal@finlytyx.com, new send date/time.
$ws.Range("A3").Value = "mashal@finlytyx.com"
$ws.Range("B3").Value = 46056
$ws.Range("B3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C3").Value = "Sent at 2026-02-03 17:09:32"
